$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (rows 2-12), keep header row formatting
$ws.Range("A2:G12").ClearContents()

# Write header row (column headers change for anaemia -> TRUE/FALSE)
$ws.Range("A1").Value2 = "profession"
$ws.Range("B1").Value2 = "FALSE_count"
$ws.Range("C1").Value2 = "FALSE_percent"
$ws.Range("D1").Value2 = "TRUE_count"
$ws.Range("E1").Value2 = "TRUE_percent"
$ws.Range("F1").Value2 = "NA_count"
$ws.Range("G1").Value2 = "NA_percent"
$ws.Range("A2").Value2 = "Hair Dresser"
$ws.Range("B2").Value2 = 2
$ws.Range("C2").Value2 = 11.11111111111111
$ws.Range("D2").Value2 = 5
$ws.Range("E2").Value2 = 6.756756756756757
$ws.Range("F2").Value2 = 2
$ws.Range("G2").Value2 = 5.128205128205128
$ws.Range("A3").Value2 = "Midwife"
$ws.Range("B3").Value2 = 1
$ws.Range("C3").Value2 = 5.555555555555555
$ws.Range("A4").Value2 = "None"
$ws.Range("B4").Value2 = 6
$ws.Range("C4").Value2 = 33.33333333333333
$ws.Range("D4").Value2 = 23
$ws.Range("E4").Value2 = 31.08108108108108
$ws.Range("F4").Value2 = 7
$ws.Range("G4").Value2 = 17.94871794871795
$ws.Range("A5").Value2 = "Seamstress"
$ws.Range("B5").Value2 = 1
$ws.Range("C5").Value2 = 5.555555555555555
$ws.Range("D5").Value2 = 4
$ws.Range("E5").Value2 = 5.405405405405405
$ws.Range("F5").Value2 = 3
$ws.Range("G5").Value2 = 7.692307692307693
$ws.Range("A6").Value2 = "Student"
$ws.Range("B6").Value2 = 1
$ws.Range("C6").Value2 = 5.555555555555555
$ws.Range("D6").Value2 = 5
$ws.Range("E6").Value2 = 6.756756756756757
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 2.564102564102564
$ws.Range("A7").Value2 = "Trader"
$ws.Range("B7").Value2 = 7
$ws.Range("C7").Value2 = 38.88888888888889
$ws.Range("D7").Value2 = 26
$ws.Range("E7").Value2 = 35.13513513513514
$ws.Range("F7").Value2 = 11
$ws.Range("G7").Value2 = 28.2051282051282
$ws.Range("A8").Value2 = "Business Owner"
$ws.Range("D8").Value2 = 1
$ws.Range("E8").Value2 = 1.351351351351351
$ws.Range("A9").Value2 = "Fishmonger"
$ws.Range("D9").Value2 = 3
$ws.Range("E9").Value2 = 4.054054054054054
$ws.Range("F9").Value2 = 6
$ws.Range("G9").Value2 = 15.38461538461539
$ws.Range("A10").Value2 = "Teacher"
$ws.Range("D10").Value2 = 3
$ws.Range("E10").Value2 = 4.054054054054054
$ws.Range("F10").Value2 = 2
$ws.Range("G10").Value2 = 5.128205128205128
$ws.Range("D11").Value2 = 4
$ws.Range("E11").Value2 = 5.405405405405405
$ws.Range("F11").Value2 = 6
$ws.Range("G11").Value2 = 15.38461538461539
$ws.Range("A12").Value2 = "Undertaker"
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 2.564102564102564
